$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2312.4
$ws.Range("I18").Value = 2544
$ws.Range("K18").Value = 2544
$ws.Range("M18").Value = -2260
$ws.Range("H40").Value = 10966.667
$ws.Range("I40").Value = 9900
$ws.Range("K40").Value = 9900
$ws.Range("M40").Value = -9725
$ws.Range("H86").Value = 2926882.2
$ws.Range("I86").Value = 4236.75
$ws.Range("J86").Value = 3761923.8
$ws.Range("K86").Value = 4236.75
$ws.Range("L86").Value = 3761923.8
$ws.Range("M86").Value = -3113.75
$ws.Range("N86").Value = -3764169.8
$ws.Range("H88").Value = 3598.6086
$ws.Range("J88").Value = 2503.8333
$ws.Range("L88").Value = 2503.8333
$ws.Range("N88").Value = -3315.8333
$ws.Range("H89").Value = 2926882.2
$ws.Range("I89").Value = 4236.75
$ws.Range("J89").Value = 3761923.8
$ws.Range("K89").Value = 21183.75
$ws.Range("L89").Value = 18809619
$ws.Range("M89").Value = -15567.75
$ws.Range("N89").Value = -18820851
$ws.Range("H91").Value = 3598.6086
$ws.Range("J91").Value = 2503.8333
$ws.Range("L91").Value = 2503.8333
$ws.Range("N91").Value = -5311.8333
$ws.Range("H131").Value = 9150.576999999999
$ws.Range("I131").Value = 7882.067
$ws.Range("K131").Value = 23646.201
$ws.Range("M131").Value = -18606.201
$ws.Range("H132").Value = 17908.795
$ws.Range("I132").Value = 1655.92
$ws.Range("J132").Value = 63055.668
$ws.Range("K132").Value = 4967.76
$ws.Range("L132").Value = 189167.004
$ws.Range("M132").Value = -2437.76
$ws.Range("N132").Value = -194227.004
$ws.Range("H137").Value = 10786.818
$ws.Range("I137").Value = 15276
$ws.Range("J137").Value = 5399.8
$ws.Range("K137").Value = 45828
$ws.Range("L137").Value = 16199.4
$ws.Range("M137").Value = -43278
$ws.Range("N137").Value = -21299.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4395.8335
$ws.Range("I61").Value = 2031.875
$ws.Range("J61").Value = 9123.75
$ws.Range("K61").Value = 2031.875
$ws.Range("L61").Value = 9123.75
$ws.Range("M61").Value = -1819.875
$ws.Range("N61").Value = -9547.75
$ws.Range("H74").Value = 3175.0476
$ws.Range("I74").Value = 3244.05
$ws.Range("J74").Value = 1795
$ws.Range("K74").Value = 3244.05
$ws.Range("L74").Value = 1795
$ws.Range("M74").Value = -2370.05
$ws.Range("N74").Value = -3543
$ws.Range("H77").Value = 3175.0476
$ws.Range("I77").Value = 3244.05
$ws.Range("J77").Value = 1795
$ws.Range("K77").Value = 16220.25
$ws.Range("L77").Value = 8975
$ws.Range("M77").Value = -11852.25
$ws.Range("N77").Value = -17711
$ws.Range("H88").Value = 4387373.5
$ws.Range("J88").Value = 6945829
$ws.Range("L88").Value = 6945829
$ws.Range("N88").Value = -6946641
$ws.Range("H91").Value = 4387373.5
$ws.Range("J91").Value = 6945829
$ws.Range("L91").Value = 6945829
$ws.Range("N91").Value = -6948637
$ws.Range("H102").Value = 16672816
$ws.Range("I102").Value = 4086.4707
$ws.Range("K102").Value = 4086.4707
$ws.Range("M102").Value = -2464.4707
$ws.Range("H136").Value = 4395.8335
$ws.Range("I136").Value = 2031.875
$ws.Range("J136").Value = 9123.75
$ws.Range("K136").Value = 6095.625
$ws.Range("L136").Value = 27371.25
$ws.Range("M136").Value = -3545.625
$ws.Range("N136").Value = -32471.25
$ws.Range("H140").Value = 76000
$ws.Range("J140").Value = 76000
$ws.Range("L140").Value = 76000
$ws.Range("N140").Value = -86360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 24310.2
$ws.Range("I96").Value = 5638
$ws.Range("K96").Value = 5638
$ws.Range("M96").Value = -2892

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10000000
$ws.Range("I6").Value = 10000000
$ws.Range("K6").Value = 10000000
$ws.Range("M6").Value = -9999887
$ws.Range("H25").Value = 1199
$ws.Range("I25").Value = 666.6667
$ws.Range("J25").Value = 1997.5
$ws.Range("K25").Value = 666.6667
$ws.Range("L25").Value = 1997.5
$ws.Range("M25").Value = -492.6667
$ws.Range("N25").Value = -2345.5
$ws.Range("H31").Value = 2876.926
$ws.Range("J31").Value = 4999.3335
$ws.Range("L31").Value = 4999.3335
$ws.Range("N31").Value = -5589.3335
$ws.Range("H34").Value = 2876.926
$ws.Range("J34").Value = 4999.3335
$ws.Range("L34").Value = 4999.3335
$ws.Range("N34").Value = -5403.3335
$ws.Range("H94").Value = 1651.3636
$ws.Range("I94").Value = 1329.6666
$ws.Range("J94").Value = 1772
$ws.Range("K94").Value = 1329.6666
$ws.Range("L94").Value = 1772
$ws.Range("M94").Value = -878.6666
$ws.Range("N94").Value = -2674

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1671949.9
$ws.Range("I9").Value = 673006.6
$ws.Range("K9").Value = 2019019.8
$ws.Range("M9").Value = -2018795.8
$ws.Range("H60").Value = 988.6667
$ws.Range("J60").Value = 3626.6667
$ws.Range("L60").Value = 10880.0001
$ws.Range("N60").Value = -11382.0001
$ws.Range("H132").Value = 1826.8572
$ws.Range("I132").Value = 1557.6
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 14018.4
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -11488.4
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 110.333336
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 182.66667
$ws.Range("K2").Value = 38
$ws.Range("L2").Value = 182.66667
$ws.Range("M2").Value = 75
$ws.Range("N2").Value = -408.66667
$ws.Range("H70").Value = 144551
$ws.Range("I70").Value = 226722
$ws.Range("J70").Value = 7599.3335
$ws.Range("K70").Value = 226722
$ws.Range("L70").Value = 7599.3335
$ws.Range("M70").Value = -226452
$ws.Range("N70").Value = -8139.3335
$ws.Range("H73").Value = 144551
$ws.Range("I73").Value = 226722
$ws.Range("J73").Value = 7599.3335
$ws.Range("K73").Value = 226722
$ws.Range("L73").Value = 7599.3335
$ws.Range("M73").Value = -225786
$ws.Range("N73").Value = -9471.333500000001
$ws.Range("H80").Value = 27874018
$ws.Range("I80").Value = 223796
$ws.Range("J80").Value = 47624176
$ws.Range("K80").Value = 223796
$ws.Range("L80").Value = 47624176
$ws.Range("M80").Value = -222798
$ws.Range("N80").Value = -47626172
$ws.Range("H83").Value = 27874018
$ws.Range("I83").Value = 223796
$ws.Range("J83").Value = 47624176
$ws.Range("K83").Value = 1118980
$ws.Range("L83").Value = 238120880
$ws.Range("M83").Value = -1113988
$ws.Range("N83").Value = -238130864
$ws.Range("H102").Value = 4617.7646
$ws.Range("I102").Value = 4150.857
$ws.Range("J102").Value = 6796.6665
$ws.Range("K102").Value = 4150.857
$ws.Range("L102").Value = 6796.6665
$ws.Range("M102").Value = -2528.857
$ws.Range("N102").Value = -10040.6665
$ws.Range("H132").Value = 10453.846
$ws.Range("I132").Value = 9933.333000000001
$ws.Range("J132").Value = 11625
$ws.Range("K132").Value = 29799.999
$ws.Range("L132").Value = 34875
$ws.Range("M132").Value = -27269.999
$ws.Range("N132").Value = -39935
$ws.Range("H134").Value = 114999.5
$ws.Range("J134").Value = 114999.5
$ws.Range("L134").Value = 344998.5
$ws.Range("N134").Value = -350068.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6062099
$ws.Range("J22").Value = 1666.0834
$ws.Range("L22").Value = 1666.0834
$ws.Range("N22").Value = -2256.0834
$ws.Range("H27").Value = 6062099
$ws.Range("J27").Value = 1666.0834
$ws.Range("L27").Value = 1666.0834
$ws.Range("N27").Value = -1880.0834
$ws.Range("H46").Value = 3169.8
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312
$ws.Range("H132").Value = 5510.625
$ws.Range("I132").Value = 4772.5
$ws.Range("K132").Value = 14317.5
$ws.Range("M132").Value = -11787.5
$ws.Range("H136").Value = 5454
$ws.Range("I136").Value = 4499.1665
$ws.Range("K136").Value = 13497.4995
$ws.Range("M136").Value = -10947.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1681.4445
$ws.Range("I122").Value = 1636.8
$ws.Range("K122").Value = 4910.4
$ws.Range("M122").Value = -2460.4
$ws.Range("H136").Value = 11000
$ws.Range("I136").Value = 11000
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 33000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -30450
$ws.Range("N136").ClearContents()
